$d = $word.ActiveDocument
$f = $d.Content
$f.Find.Execute("ootnotefay")
Write-Host "ootnotefay: start=$($f.Start) end=$($f.End)"
for ($i=720; $i -le 726; $i++) {
  $r = $d.Range($i, $i+1)
  Write-Host "$i : [$($r.Text)] ($($r.End - $r.Start))"
}
